# edit.ps1 - reproduces the authored change:
#   1. slide6's table switches from the deck's custom "Table_0" style
#      ({D45E5890-CF39-488E-9B14-6F441AA06E3B}) to the built-in table
#      style {904916DF-CB55-4AFE-9215-670F7683374F}.
#   2. The presentation's theme palette (ppt/theme/theme2.xml, the theme
#      part actually referenced by the slide master / presentation) is
#      switched from the "Integral" palette to the standard "Office
#      Theme" palette (the palette that used to live, unused by any
#      slide, in ppt/theme/theme1.xml).

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{904916DF-CB55-4AFE-9215-670F7683374F}")
    }
}

# --- 2. Theme colors --------------------------------------------------
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$cs.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$cs.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$cs.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$cs.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$cs.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$cs.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$cs.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$cs.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$cs.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$cs.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$cs.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
